$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The diff adds one new data row ("Coliflor" observation dated 45015 /
# 2023-03-30) at the top of the data block that starts on row 417. Every
# existing row from 417 to 507 is pushed down by one (417->418, ..., 507->508),
# and the sheet's used range grows from A1:R507 to A1:R508.

# 1) Insert a blank row at 417; rows 417:507 shift down to 418:508.
$ws.Rows.Item(417).Insert()

# 2) The row that landed on 418 (the former row 417) has the same
#    look/format/values as the new row we need on 417 except for the
#    "Volumen" (J) and "Fecha" (D) fields, so copy it up to seed row 417
#    with correct styles/shared values, then overwrite D/J with the new data.
$ws.Rows.Item(418).Copy()
$ws.Rows.Item(417).PasteSpecial()

$ws.Cells.Item(417, 4).Value = 45015
$ws.Cells.Item(417, 10).Value = 500
